# The "override" column (E) was populated for every row as a holdover
# from a previous merge step. Per the commit message ("merge sheets are
# not run if any overrides are there"), the override values for the data
# rows should be cleared so the merge step isn't skipped. The header in
# E1 ("override") is left untouched; only the data cells E2:E230 are
# cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) {
    $lastRow = 230
}

$ws.Range("E2:E$lastRow").ClearContents()
